$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4953129999999999
$ws.Range("H2").Value = 1.485939
$ws.Range("I2").Value = 0.2134181161789063
$ws.Range("J2").Value = 0.2134181161789063
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.345785
$ws.Range("N2").Value = 1.037355
$ws.Range("O2").Value = 0.1901295499320662
$ws.Range("P2").Value = 0.1901295499320662
$ws.Range("Q2").Value = 0.171271805705
$ws.Range("R2").Value = 1.541446251345
$ws.Range("S2").Value = 0.04057709037644486
$ws.Range("T2").Value = 0.04057709037644486

# Row 3
$ws.Range("G3").Value = 0.4953129999999999
$ws.Range("H3").Value = 1.485939
$ws.Range("I3").Value = 0.2134181161789063
$ws.Range("J3").Value = 0.2134181161789063
$ws.Range("O3").Value = 0.04212778381695306
$ws.Range("P3").Value = 0.04212778381695306
$ws.Range("Q3").Value = 0.037949396121
$ws.Range("R3").Value = 0.341544565089
$ws.Range("S3").Value = 0.008990832261006337
$ws.Range("T3").Value = 0.008990832261006337

# Row 4
$ws.Range("G4").Value = 0.4953129999999999
$ws.Range("H4").Value = 1.485939
$ws.Range("I4").Value = 0.2134181161789063
$ws.Range("J4").Value = 0.2134181161789063
$ws.Range("M4").Value = 1.396279
$ws.Range("N4").Value = 4.188836999999999
$ws.Range("O4").Value = 0.7677426662509808
$ws.Range("P4").Value = 0.7677426662509808
$ws.Range("Q4").Value = 0.6915951403269999
$ws.Range("R4").Value = 6.224356262942999
$ws.Range("S4").Value = 0.1638501935414551
$ws.Range("T4").Value = 0.1638501935414551

# Row 5
$ws.Range("I5").Value = 0.3878984662564351
$ws.Range("J5").Value = 0.3878984662564351
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.345785
$ws.Range("N5").Value = 1.037355
$ws.Range("O5").Value = 0.1901295499320662
$ws.Range("P5").Value = 0.1901295499320662
$ws.Range("Q5").Value = 0.311295366745
$ws.Range("R5").Value = 2.801658300705
$ws.Range("S5").Value = 0.07375096080867477
$ws.Range("T5").Value = 0.07375096080867477

# Row 6
$ws.Range("I6").Value = 0.3878984662564351
$ws.Range("J6").Value = 0.3878984662564351
$ws.Range("O6").Value = 0.04212778381695306
$ws.Range("P6").Value = 0.04212778381695306
$ws.Range("R6").Value = 0.620774915121
$ws.Range("S6").Value = 0.01634130272937876
$ws.Range("T6").Value = 0.01634130272937876

# Row 7
$ws.Range("I7").Value = 0.3878984662564351
$ws.Range("J7").Value = 0.3878984662564351
$ws.Range("M7").Value = 1.396279
$ws.Range("N7").Value = 4.188836999999999
$ws.Range("O7").Value = 0.7677426662509808
$ws.Range("P7").Value = 0.7677426662509808
$ws.Range("Q7").Value = 1.257009943703
$ws.Range("R7").Value = 11.313089493327
$ws.Range("S7").Value = 0.2978062027183815
$ws.Range("T7").Value = 0.2978062027183816

# Row 8
$ws.Range("G8").Value = 0.8083133333333334
$ws.Range("H8").Value = 2.42494
$ws.Range("I8").Value = 0.3482822152503414
$ws.Range("J8").Value = 0.3482822152503415
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.345785
$ws.Range("N8").Value = 1.037355
$ws.Range("O8").Value = 0.1901295499320662
$ws.Range("P8").Value = 0.1901295499320662
$ws.Range("Q8").Value = 0.2795026259666667
$ws.Range("R8").Value = 2.5155236337
$ws.Range("S8").Value = 0.0662187408348904
$ws.Range("T8").Value = 0.06621874083489042

# Row 9
$ws.Range("G9").Value = 0.8083133333333334
$ws.Range("H9").Value = 2.42494
$ws.Range("I9").Value = 0.3482822152503414
$ws.Range("J9").Value = 0.3482822152503415
$ws.Range("O9").Value = 0.04212778381695306
$ws.Range("P9").Value = 0.04212778381695306
$ws.Range("Q9").Value = 0.06193054266000001
$ws.Range("R9").Value = 0.5573748839400001
$ws.Range("S9").Value = 0.0146723578713559
$ws.Range("T9").Value = 0.0146723578713559

# Row 10
$ws.Range("G10").Value = 0.8083133333333334
$ws.Range("H10").Value = 2.42494
$ws.Range("I10").Value = 0.3482822152503414
$ws.Range("J10").Value = 0.3482822152503415
$ws.Range("M10").Value = 1.396279
$ws.Range("N10").Value = 4.188836999999999
$ws.Range("O10").Value = 0.7677426662509808
$ws.Range("P10").Value = 0.7677426662509808
$ws.Range("Q10").Value = 1.128630932753333
$ws.Range("R10").Value = 10.15767839478
$ws.Range("S10").Value = 0.2673911165440951
$ws.Range("T10").Value = 0.2673911165440951

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.116974
$ws.Range("H11").Value = 0.350922
$ws.Range("I11").Value = 0.05040120231431718
$ws.Range("J11").Value = 0.05040120231431718
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.345785
$ws.Range("N11").Value = 1.037355
$ws.Range("O11").Value = 0.1901295499320662
$ws.Range("P11").Value = 0.1901295499320662
$ws.Range("Q11").Value = 0.04044785459000001
$ws.Range("R11").Value = 0.36403069131
$ws.Range("S11").Value = 0.009582757912056138
$ws.Range("T11").Value = 0.009582757912056138

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.116974
$ws.Range("H12").Value = 0.350922
$ws.Range("I12").Value = 0.05040120231431718
$ws.Range("J12").Value = 0.05040120231431718
$ws.Range("O12").Value = 0.04212778381695306
$ws.Range("P12").Value = 0.04212778381695306
$ws.Range("Q12").Value = 0.008962196958000001
$ws.Range("R12").Value = 0.080659772622
$ws.Range("S12").Value = 0.002123290955212069
$ws.Range("T12").Value = 0.002123290955212069

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.116974
$ws.Range("H13").Value = 0.350922
$ws.Range("I13").Value = 0.05040120231431718
$ws.Range("J13").Value = 0.05040120231431718
$ws.Range("M13").Value = 1.396279
$ws.Range("N13").Value = 4.188836999999999
$ws.Range("O13").Value = 0.7677426662509808
$ws.Range("P13").Value = 0.7677426662509808
$ws.Range("Q13").Value = 0.163328339746
$ws.Range("R13").Value = 1.469955057714
$ws.Range("S13").Value = 0.03869515344704898
$ws.Range("T13").Value = 0.03869515344704898
